# Update the session-results table:
#  - replace the three student names with a new roster of six names
#  - rewrite the Group/Exam/Session/Mark figures for rows 2-7
#  - drop the old rows 8 & 9 (two extra records no longer present)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Karp Alex DD"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 3

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "NIkita Valer Gnusov"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "MOyva JJu DD"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 10

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Prerk Tuk Tuc"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 8

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Munsd sdhf dfkj"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 8

# Row 7
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Larisa JJJ Ddsf"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2

# The old row 9 (Tsmyg Dmitry Alexandrovich, session 5) is removed entirely,
# shifting nothing below it up (it's the last row); row 8 (Tsmyg Dmitry
# Alexandrovich, session 4) is cleared out so the table ends one row sooner.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).ClearContents()
